$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Fill in row 11 with the new time-tracking entry
$ws.Range("E11").Value = "GL3"
$ws.Range("F11").Value = "GDPR"
$ws.Range("G11").Value = "Riunione interna"
$ws.Range("H11").Value = (Get-Date -Year 2019 -Month 3 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I11").Value = 90

# Move the active selection to J11, matching the saved cursor position
$ws.Range("J11").Select()
